$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 11 with the new contribution log entry
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "Ryan Conyac"
$ws.Range("D11").Value = "Edited video for sprint 2"

# Update the active selection to D12, matching the recorded sheet view
$ws.Range("D12").Select()
